$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 28: pyruvate synthase (PFOR) gene entry, appended after PSBS2
$ws.Range("A28").Value = "Cre06.g292250"
$ws.Range("C28").Value = "PFOR"
$ws.Range("D28").Value = "pyruvate synthase"
$ws.Range("G28").Value = "Hydrogen_Prod"

# Copy the existing table's cell formatting onto the new row so the added
# cells match the look of the rest of the sheet (center-aligned id/symbol
# cells, left+vcenter pathway cell) instead of picking up default formatting
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)

$ws.Range("C27").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("G27").Copy()
$ws.Range("G28").PasteSpecial(-4122)

$ws.Range("A27").Copy()
$ws.Range("A29").PasteSpecial(-4122)

# Widen column G a bit so the longer pathway label fits
$ws.Columns.Item(7).ColumnWidth = 14.6666666666667

# Leave the selection where the author ended up after entering the data
$ws.Range("C32").Select()
